# "Se procesan de nuevo los datos con las nuevas dimensiones curadas"
#
# The curated-dimensions re-processing changes three columns of the
# metadata header table (rows 2-5 describe each data column in row 1):
#
#   - Column H (municipio-nombre): was wired up as an
#     "iaest-measure:municipio-nombre" measure; now it is curated as a
#     proper dimension, matching column I (provincia-nombre):
#       H2 medida-uri -> sdmx-dimension:refArea
#       H3 medida/dim -> dim
#       H4 type       -> URI-Municipio
#
#   - Column U (total-bienes-inmuebles): was curated as a dimension
#     (with a concept mapping file); now it is a plain measure:
#       U2 iaest-dimension:total-bienes-inmuebles -> iaest-measure:total-bienes-inmuebles
#       U3 dim -> medida
#       U4 skos:Concept -> xsd:int
#       U5 mapping-total-bienes-inmuebles.xlsx -> (cleared)
#
#   - Column X (bi-espectaculos): same kind of dimension -> measure fix:
#       X2 iaest-dimension:bi-espectaculos -> iaest-measure:bi-espectaculos
#       X3 dim -> medida
#       X4 skos:Concept -> xsd:int
#       X5 mapping-bi-espectaculos.xlsx -> (cleared)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H (municipio-nombre) becomes a curated dimension like column I.
$ws.Range("H2").Value = "sdmx-dimension:refArea"
$ws.Range("H3").Value = "dim"
$ws.Range("H4").Value = "URI-Municipio"

# Column U (total-bienes-inmuebles) becomes a plain measure.
$ws.Range("U2").Value = "iaest-measure:total-bienes-inmuebles"
$ws.Range("U3").Value = "medida"
$ws.Range("U4").Value = "xsd:int"
$ws.Range("U5").ClearContents()

# Column X (bi-espectaculos) becomes a plain measure.
$ws.Range("X2").Value = "iaest-measure:bi-espectaculos"
$ws.Range("X3").Value = "medida"
$ws.Range("X4").Value = "xsd:int"
$ws.Range("X5").ClearContents()
